$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("06-09-2021", 4.86, 4.17, 3.56, 3.38),
    @("07-09-2021", 4.73, 4.13, 3.56, 3.39),
    @("08-09-2021", 4.81, 4.11, 3.59, 3.43),
    @("09-09-2021", 4.88, 4.11, 3.63, 3.46),
    @("10-09-2021", 4.81, 4.07, 3.65, 3.47)
)

$startRow = 175
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $data[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
